$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'332.40"
$ws.Range("E2").Value = "'1.06%"
$ws.Range("D3").Value = "'45.61"
$ws.Range("E3").Value = "'3.07%"
$ws.Range("D4").Value = "'5.473"
$ws.Range("E4").Value = "'-0.53%"
$ws.Range("D5").Value = "'0.08523"
$ws.Range("E5").Value = "'5.54%"
$ws.Range("D6").Value = "'2.057"
$ws.Range("E6").Value = "'0.51%"
$ws.Range("D7").Value = "'0.9875"
$ws.Range("E7").Value = "'3.66%"
$ws.Range("D9").Value = "'0.1152"
$ws.Range("E9").Value = "'0.19%"
$ws.Range("D10").Value = "'0.1913"
$ws.Range("E10").Value = "'1.70%"
$ws.Range("D11").Value = "'9.484"
$ws.Range("E11").Value = "'-6.91%"
$ws.Range("D12").Value = "'0.09749"
$ws.Range("E12").Value = "'-2.01%"
$ws.Range("D13").Value = "'0.04723"
$ws.Range("E13").Value = "'-2.30%"
$ws.Range("E14").Value = "'-0.36%"
$ws.Range("D15").Value = "'0.001302"
$ws.Range("E15").Value = "'3.78%"
$ws.Range("D16").Value = "'0.005986"
$ws.Range("E16").Value = "'1.93%"
$ws.Range("D17").Value = "'3.384"
$ws.Range("E17").Value = "'0.28%"
$ws.Range("D18").Value = "'4.448"
$ws.Range("E18").Value = "'1.02%"
$ws.Range("D19").Value = "'0.3354"
$ws.Range("E19").Value = "'-1.50%"
$ws.Range("D20").Value = "'0.1374"
$ws.Range("E20").Value = "'-1.91%"
$ws.Range("E21").Value = "'-1.02%"
$ws.Range("D22").Value = "'0.04143"
$ws.Range("E22").Value = "'1.42%"
$ws.Range("E23").Value = "'-0.29%"
$ws.Range("D24").Value = "'0.004485"
$ws.Range("E24").Value = "'3.05%"
$ws.Range("D25").Value = "'0.0001303"
$ws.Range("E25").Value = "'4.24%"
$ws.Range("E26").Value = "'-20.11%"
$ws.Range("D38").Value = "'0.02769"
$ws.Range("E38").Value = "'6.76%"
$ws.Range("D39").Value = "'0.05707"
$ws.Range("E39").Value = "'0.06%"
$ws.Range("D40").Value = "'0.007906"
$ws.Range("E40").Value = "'4.37%"
$ws.Range("D41").Value = "'0.1432"
$ws.Range("E41").Value = "'2.39%"
$ws.Range("D42").Value = "'0.007264"
$ws.Range("E42").Value = "'-0.93%"
$ws.Range("D43").Value = "'0.002113"
$ws.Range("E43").Value = "'5.27%"
$ws.Range("D44").Value = "'0.009003"
$ws.Range("E44").Value = "'-0.46%"
$ws.Range("D45").Value = "'0.3554"
$ws.Range("D46").Value = "'0.00006993"
$ws.Range("E46").Value = "'0.02%"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("E47").Value = "'0.23%"
$ws.Range("E48").Value = "'0.31%"
$ws.Range("D49").Value = "'0.003459"
$ws.Range("E49").Value = "'-0.92%"
$ws.Range("D50").Value = "'0.003537"
$ws.Range("E50").Value = "'1.14%"
$ws.Range("E51").Value = "'0.23%"
